# Actualización automática desde tarea programada
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 / column A: refreshed timestamp (tiny recalculated precision change)
$ws.Cells.Item(5, 1).Value2 = 45869.62521966435

# New row 6: next scheduled sensor reading
$ws.Cells.Item(6, 1).Value2 = 45869.75026166721
$ws.Cells.Item(6, 1).NumberFormat = $ws.Cells.Item(5, 1).NumberFormat

$ws.Cells.Item(6, 2).Value2 = 2025
$ws.Cells.Item(6, 3).Value2 = 31
$ws.Cells.Item(6, 4).Value2 = 14.94
$ws.Cells.Item(6, 5).Value2 = 88.26000000000001
$ws.Cells.Item(6, 6).Value2 = 3.52
$ws.Cells.Item(6, 7).Value2 = 10.85
$ws.Cells.Item(6, 8).Value2 = "ESE"
$ws.Cells.Item(6, 9).Value2 = 0
$ws.Cells.Item(6, 10).Value2 = "18:00:22"
